$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, pushing the existing rows 86-101 down to 87-102.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly record.
$ws.Cells.Item(86, 1).Value()  = 1
$ws.Cells.Item(86, 2).Value()  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(86, 3).Value()  = "Arica y Parinacota"
$ws.Cells.Item(86, 4).Value()  = 44816
$ws.Cells.Item(86, 5).Value()  = 15
$ws.Cells.Item(86, 6).Value()  = 100112038
$ws.Cells.Item(86, 7).Value()  = "Cebollín baby"
$ws.Cells.Item(86, 8).Value()  = "Sin especificar"
$ws.Cells.Item(86, 9).Value()  = "Primera"
$ws.Cells.Item(86, 10).Value() = 250
$ws.Cells.Item(86, 11).Value() = 1200
$ws.Cells.Item(86, 12).Value() = 1500
$ws.Cells.Item(86, 13).Value() = 1350
$ws.Cells.Item(86, 14).Value() = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(86, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value() = 675
$ws.Cells.Item(86, 17).Value() = 2
$ws.Cells.Item(86, 18).Value() = "Hortaliza"
